$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''59.386.43'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '''  -0.30%  '
$ws.Range("E2").Style = "Normal"

$ws.Range("D3").Value = '''2.517.99'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '''  +1.08%  '
$ws.Range("E3").Style = "Normal"

$ws.Range("E4").Value = '''  -0.18%  '
$ws.Range("E4").Style = "Normal"

$ws.Range("D5").Value = '''542.80'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '''  +0.02%  '
$ws.Range("E5").Style = "Normal"

$ws.Range("D6").Value = '''145.17'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '''  -1.15%  '
$ws.Range("E6").Style = "Normal"

$ws.Range("E7").Value = '''  -0.20%  '
$ws.Range("E7").Style = "Normal"

$ws.Range("D8").Value = '''0.573'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '''  -0.88%  '
$ws.Range("E8").Style = "Normal"

$ws.Range("D9").Value = '''2.539.03'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '''  +0.87%  '
$ws.Range("E9").Style = "Normal"

$ws.Range("E10").Value = '''  -0.25%  '
$ws.Range("E10").Style = "Normal"

$ws.Range("E11").Value = '''  +0.19%  '
$ws.Range("E11").Style = "Normal"

$ws.Range("D12").Value = '''5.58'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '''  +1.72%  '
$ws.Range("E12").Style = "Normal"

$ws.Range("D13").Value = '''0.361'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '''  +1.80%  '
$ws.Range("E13").Style = "Normal"

$ws.Range("D14").Value = '''2.962.35'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '''  -0.21%  '
$ws.Range("E14").Style = "Normal"

$ws.Range("D15").Value = '''23.51'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '''  -3.95%  '
$ws.Range("E15").Style = "Normal"

$ws.Range("D16").Value = '''59.278.31'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '''  -0.65%  '
$ws.Range("E16").Style = "Normal"

$ws.Range("E17").Value = '''  +1.22%  '
$ws.Range("E17").Style = "Normal"

$ws.Range("D18").Value = '''2.534.10'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '''  +1.40%  '
$ws.Range("E18").Style = "Normal"

$ws.Range("E19").Value = '''  -1.76%  '
$ws.Range("E19").Style = "Normal"

$ws.Range("D21").Value = '''326.50'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '''  -0.06%  '
$ws.Range("E21").Style = "Normal"

$ws.Range("E22").Value = '''  +0.38%  '
$ws.Range("E22").Style = "Normal"

$ws.Range("E23").Value = '''  +1.16%  '
$ws.Range("E23").Style = "Normal"

$ws.Range("D24").Value = '''62.19'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '''  +1.35%  '
$ws.Range("E24").Style = "Normal"

$ws.Range("D25").Value = '''0.433'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '''  -2.86%  '
$ws.Range("E25").Style = "Normal"

$ws.Range("E26").Value = '''  +1.63%  '
$ws.Range("E26").Style = "Normal"

$ws.Range("E27").Value = '''  -1.42%  '
$ws.Range("E27").Style = "Normal"

$ws.Range("D28").Value = '''8.02'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '''  +2.44%  '
$ws.Range("E28").Style = "Normal"

$ws.Range("D29").Value = '''0.0₃0787'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '''  -1.09%  '
$ws.Range("E29").Style = "Normal"

$ws.Range("D30").Value = '''1.82'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '''  -0.07%  '
$ws.Range("E30").Style = "Normal"

$ws.Range("D31").Value = '''6.73'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '''  -0.73%  '
$ws.Range("E31").Style = "Normal"

$ws.Range("D32").Value = '''1.20'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '''  -6.55%  '
$ws.Range("E32").Style = "Normal"

$ws.Range("D33").Value = '''1.48'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '''  +2.51%  '
$ws.Range("E33").Style = "Normal"

$ws.Range("E34").Value = '''  +0.09%  '
$ws.Range("E34").Style = "Normal"

$ws.Range("D35").Value = '''159.27'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '''  +0.75%  '
$ws.Range("E35").Style = "Normal"

$ws.Range("D36").Value = '''18.76'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '''  -1.59%  '
$ws.Range("E36").Style = "Normal"

$ws.Range("E37").Value = '''  -1.98%  '
$ws.Range("E37").Style = "Normal"

$ws.Range("E38").Value = '''  -7.04%  '
$ws.Range("E38").Style = "Normal"

$ws.Range("D39").Value = '''37.05'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '''  +0.80%  '
$ws.Range("E39").Style = "Normal"

$ws.Range("D40").Value = '''5.58'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '''  -5.95%  '
$ws.Range("E40").Style = "Normal"

$ws.Range("D41").Value = '''0.835'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '''  +0.65%  '
$ws.Range("E41").Style = "Normal"

$ws.Range("D42").Value = '''295.54'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '''  -6.23%  '
$ws.Range("E42").Style = "Normal"

$ws.Range("D43").Value = '''3.69'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '''  -2.38%  '
$ws.Range("E43").Style = "Normal"

$ws.Range("E44").Value = '''  -0.06%  '
$ws.Range("E44").Style = "Normal"

$ws.Range("E45").Value = '''  -0.32%  '
$ws.Range("E45").Style = "Normal"

$ws.Range("E46").Value = '''  +1.28%  '
$ws.Range("E46").Style = "Normal"

$ws.Range("D47").Value = '''0.0937'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '''  -0.42%  '
$ws.Range("E47").Style = "Normal"

$ws.Range("D48").Value = '''18.85'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '''  +0.67%  '
$ws.Range("E48").Style = "Normal"

$ws.Range("D49").Value = '''123.19'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '''  -3.05%  '
$ws.Range("E49").Style = "Normal"

$ws.Range("E50").Value = '''  -2.21%  '
$ws.Range("E50").Style = "Normal"

$ws.Range("E51").Value = '''  -3.66%  '
$ws.Range("E51").Style = "Normal"
